{"js": "const body = context.document.body;\n\nconst title = \"H\u01b0\u1edbng d\u1eabn c\u00e0i \u0111\u1eb7t WORDPRESS tr\u00ean SUBDOMAIN\";\nconst body2 = \"SUBDOMAIN kh\u00f4ng h\u01b0\u1edbng d\u1eabn do thu\u1ed9c v\u1ec1 wordpress.com\";\n\n// Wipe the whole body first - this also drops the stray <w:proofErr/>\n// spell-check markers and the old \"_GoBack\" bookmark that were attached\n// to the original single paragraph.\nbody.clear();\n\nconst p0 = body.paragraphs.getFirst();\np0.insertText(title, \"Replace\");\n\n// Second paragraph with the replacement text.\nconst p1 = p0.insertParagraph(body2, \"After\");\nawait context.sync();\n\n// Recreate the \"_GoBack\" bookmark wrapping paragraph 2's text (matches the\n// position it held relative to the edited text in the original document).\nconst p1Range = p1.getRange();\np1Range.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$title = \"H\u01b0\u1edbng d\u1eabn c\u00e0i \u0111\u1eb7t WORDPRESS tr\u00ean SUBDOMAIN\"\n$body2 = \"SUBDOMAIN kh\u00f4ng h\u01b0\u1edbng d\u1eabn do thu\u1ed9c v\u1ec1 wordpress.com\"\n\n# Wipe the whole document body (this also drops the stray <w:proofErr/>\n# spell-check markers and the old \"_GoBack\" bookmark that were attached\n# to the original single paragraph).\n$r = $d.Content\n$r.Delete()\n\n# Paragraph 1: plain title text, no proof-err markers.\n$r.InsertAfter($title)\n\n# Start paragraph 2.\n$r.InsertParagraphAfter()\n$r.Collapse(0)\n\n# Paragraph 2: replacement text.\n$r.InsertAfter($body2)\n\n# Recreate the \"_GoBack\" bookmark wrapping the paragraph-2 text (matches\n# the position it held relative to the edited text in the original doc).\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
